$wb = $excel.ActiveWorkbook

# --- Rename the original sheet, add the new "Misc Sizing" sheet right after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Rocket Sizing"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Misc Sizing"

# --- Pre-format the area-ratio matrix region so blank cells still carry the
#     "0.00" number style, matching how the grid was built in Excel (format the
#     range first, then fill in only the upper-triangle formulas/values). ---
$ws2.Range("C3:G7").NumberFormat = "0.00"
$ws2.Range("C2:G2").NumberFormat = "0.00"
$ws2.Range("C2:G2").Font.Bold = $true
$ws2.Range("B3:B7").Font.Bold = $true

# --- Title ---
$ws2.Range("B2").Value = "ID\OD"
$ws2.Range("A1").Value = "Area Ratios in Concentric Tubes"
$ws2.Range("A1").Font.Bold = $true

# --- Column/row headers of the matrix (ID/OD sizes) ---
$ws2.Range("C2").Value = 0.25
$ws2.Range("D2").Formula = "=3/8"
$ws2.Range("E2").Value = 0.5
$ws2.Range("F2").Formula = "=3/4"
$ws2.Range("G2").Value = 1

$ws2.Range("B3").Value = 0.25
$ws2.Range("B4").Formula = "=3/8"
$ws2.Range("B5").Value = 0.5
$ws2.Range("B6").Formula = "=3/4"
$ws2.Range("B7").Value = 1

# --- Upper-triangle area-ratio formulas, filled right then down a step at a time ---
$ws2.Range("C3").Formula = "=((C2/2)^2-(`$B3/2)^2)/(`$B3/2)^2"
$ws2.Range("D3:G3").Formula = "=((D2/2)^2-(`$B3/2)^2)/(`$B3/2)^2"

$ws2.Range("D4").Formula = "=((D2/2)^2-(`$B4/2)^2)/(`$B4/2)^2"
$ws2.Range("E4:G4").Formula = "=((E2/2)^2-(`$B4/2)^2)/(`$B4/2)^2"

$ws2.Range("E5").Formula = "=((E2/2)^2-(`$B5/2)^2)/(`$B5/2)^2"
$ws2.Range("F5:G5").Formula = "=((F2/2)^2-(`$B5/2)^2)/(`$B5/2)^2"

$ws2.Range("F6").Formula = "=((F2/2)^2-(`$B6/2)^2)/(`$B6/2)^2"
$ws2.Range("G6").Formula = "=((G2/2)^2-(`$B6/2)^2)/(`$B6/2)^2"

$ws2.Range("G7").Value = 0

# --- Burst disk design notes ---
$ws2.Range("A10").Value = "Burst Disk Design"

$ws2.Range("A11").Value = "ORFS ID"
$ws2.Range("C11").Value = "0.437500in"
$ws2.Range("A12").Value = "Face OD"
$ws2.Range("C12").Value = "0.718750in"
$ws2.Range("A13").Value = "ORFS Thread"
$ws2.Range("A14").Value = "ORB Thread"
$ws2.Range("C14").Value = '3/4"-16'
$ws2.Range("C14").NumberFormat = "General"
$ws2.Range("C13").Value = '13/16"-16'
$ws2.Range("B11").Value = "in"
$ws2.Range("B12").Value = "in"

# --- Column width (matches the best-fit width Excel computed for column A) / view bits ---
$ws2.Columns("A").ColumnWidth = 10.14

$ws2.Range("C13").Select() | Out-Null
